$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Converted Data")

# --- Row 7 updates: U7 1 -> 0, W7 13 -> 12 ---
$ws.Cells.Item(7, 21).Value = 0
$ws.Cells.Item(7, 23).Value = 12

# --- W column (col 23) updates for rows 24-221 ---
$wValues = @{
    24 = 0.04166666666666666
    25 = 0.04166666666666666
    26 = 0.04166666666666666
    27 = 0.125
    28 = 0.125
    29 = 0.125
    30 = 0.125
    31 = 0.125
    32 = 0.2916666666666667
    33 = 0.2916666666666667
    34 = 0.2916666666666667
    35 = 0.2916666666666667
    36 = 0.2916666666666667
    37 = 0.2916666666666667
    38 = 0.2916666666666667
    39 = 0.2916666666666667
    40 = 0.2916666666666667
    41 = 0.2916666666666667
    42 = 0.8333333333333334
    43 = 0.8333333333333334
    44 = 0.8333333333333334
    45 = 0.8333333333333334
    46 = 0.8333333333333334
    47 = 0.8333333333333334
    48 = 0.8333333333333334
    49 = 0.8333333333333334
    50 = 0.8333333333333334
    51 = 0.8333333333333334
    52 = 0.8333333333333334
    53 = 0.8333333333333334
    54 = 0.8333333333333334
    55 = 0.8333333333333334
    56 = 0.8333333333333334
    57 = 0.8333333333333334
    58 = 0.8333333333333334
    59 = 0.8333333333333334
    60 = 0.8333333333333334
    61 = 0.8333333333333334
    62 = 0.8333333333333334
    63 = 0.8333333333333334
    64 = 0.8333333333333334
    65 = 0.8333333333333334
    66 = 0.8333333333333334
    67 = 0.6666666666666666
    68 = 0.6666666666666666
    69 = 0.6666666666666666
    70 = 0.6666666666666666
    71 = 0.6666666666666666
    72 = 0.6666666666666666
    73 = 0.6666666666666666
    74 = 0.6666666666666666
    75 = 0.6666666666666666
    76 = 0.6666666666666666
    77 = 0.5
    78 = 0.5
    79 = 0.5
    80 = 0.5
    81 = 0.3611111111083334
    82 = 0.3611111111083334
    83 = 0.3611111111083334
    84 = 0.3611111111083334
    85 = 0.3333333333333333
    86 = 0.3333333333333333
    87 = 0.3333333333333333
    88 = 0.3333333333333333
    89 = 0.3333333333333333
    90 = 0.3333333333333333
    91 = 0.2916666666666667
    92 = 0.2916666666666667
    93 = 0.2916666666666667
    94 = 0.2916666666666667
    95 = 0.2916666666666667
    96 = 0.2916666666666667
    97 = 0.2916666666666667
    98 = 0.2916666666666667
    99 = 0.2916666666666667
    100 = 0.2916666666666667
    101 = 0.2916666666666667
    102 = 0.25
    103 = 0.25
    104 = 0.25
    105 = 0.25
    106 = 0.25
    107 = 0.25
    108 = 0.25
    109 = 0.25
    110 = 0.25
    111 = 0.25
    112 = 0.25
    113 = 0.25
    114 = 0.25
    115 = 0.25
    116 = 0.25
    117 = 0.25
    118 = 0.25
    119 = 0.25
    120 = 0.25
    121 = 0.25
    122 = 0.25
    123 = 0.25
    124 = 0.25
    125 = 0.25
    126 = 0.25
    127 = 0.25
    128 = 0.25
    129 = 0.25
    130 = 0.25
    131 = 0.25
    132 = 0.25
    133 = 0.25
    134 = 0.25
    135 = 0.25
    136 = 0.25
    137 = 0.25
    138 = 0.25
    139 = 0.25
    140 = 0.25
    141 = 0.25
    142 = 0.25
    143 = 0.25
    144 = 0.25
    145 = 0.25
    146 = 0.25
    147 = 0.25
    148 = 0.25
    149 = 0.25
    150 = 0.25
    151 = 0.25
    152 = 0.25
    153 = 0.25
    154 = 0.25
    155 = 0.25
    156 = 0.25
    157 = 0.25
    158 = 0.25
    159 = 0.25
    160 = 0.25
    161 = 0.25
    162 = 0.25
    163 = 0.25
    164 = 0.25
    165 = 0.25
    166 = 0.25
    167 = 0.25
    168 = 0.25
    169 = 0.25
    170 = 0.25
    171 = 0.25
    172 = 0.25
    173 = 0.25
    174 = 0.25
    175 = 0.25
    176 = 0.25
    177 = 0.25
    178 = 0.25
    179 = 0.25
    180 = 0.25
    181 = 0.25
    182 = 0.25
    183 = 0.25
    184 = 0.25
    185 = 0.25
    186 = 0.25
    187 = 0.25
    188 = 0.25
    189 = 0.25
    190 = 0.25
    191 = 0.25
    192 = 0.25
    193 = 0.25
    194 = 0.25
    195 = 0.25
    196 = 0.25
    197 = 0.25
    198 = 0.25
    199 = 0.25
    200 = 0.25
    201 = 0.25
    202 = 0.25
    203 = 0.25
    204 = 0.25
    205 = 0.25
    206 = 0.25
    207 = 0.25
    208 = 0.25
    209 = 0.25
    210 = 0.25
    211 = 0.25
    212 = 0.25
    213 = 0.25
    214 = 0.25
    215 = 0.25
    216 = 0.25
    217 = 0.25
    218 = 0.25
    219 = 0.25
    220 = 0.25
    221 = 0.25
}
foreach ($r in $wValues.Keys) {
    $ws.Cells.Item($r, 23).Value = $wValues[$r]
}

# --- New rows 222-233 ---
$newRowVals = @(0,0,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,0)
$newRowDates = @(
    "9/30/2020",
    "10/1/2020",
    "10/2/2020",
    "10/3/2020",
    "10/4/2020",
    "10/5/2020",
    "10/6/2020",
    "10/7/2020",
    "10/8/2020",
    "10/9/2020",
    "10/10/2020",
    "10/11/2020",
)

$startRow = 222
for ($i = 0; $i -lt $newRowDates.Count; $i++) {
    $r = $startRow + $i
    $dateStr = $newRowDates[$i]

    # Set date text in column A, forcing text type, then restore formatting from A221
    $ws.Cells.Item($r, 1).Value = "'" + $dateStr
    $ws.Range("A221").Copy()
    $ws.Cells.Item($r, 1).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

    # Columns B..V (2..22)
    for ($c = 2; $c -le 22; $c++) {
        $ws.Cells.Item($r, $c).Value = $newRowVals[$c - 2]
    }
    # Column W (23) = 0.25
    $ws.Cells.Item($r, 23).Value = 0.25
}

$excel.CutCopyMode = $false
Write-Host "Done"
